$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.163.79"
$ws.Range("E2").Value = "  -6.24%  "
$ws.Range("D3").Value = "2.443.32"
$ws.Range("E3").Value = "  -9.03%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "532.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.83%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.567"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.89%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0987"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.63%  "
$ws.Range("E10").Value = "  -2.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.36"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.349"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.46%  "
$ws.Range("D13").Value = "2.876.77"
$ws.Range("E13").Value = "  -8.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "24.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -8.09%  "
$ws.Range("D15").Value = "59.120.58"
$ws.Range("E15").Value = "  -6.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000138"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.76%  "
$ws.Range("D17").Value = "2.491.21"
$ws.Range("E17").Value = "  -7.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -7.21%  "
$ws.Range("E19").Value = "  -5.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "323.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.84%  "
$ws.Range("E21").Value = "  -3.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.66"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -10.56%  "
$ws.Range("E23").Value = "  -7.64%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E25").Value = "  -4.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.975"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.31%  "
$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.09%  "
$ws.Range("B29").Value = "Fetch.AI"
$ws.Range("C29").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.82"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.82%  "
$ws.Range("D31").Value = "0.0₃0771"
$ws.Range("E31").Value = "  -10.44%  "
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "157.02"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.63%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.53"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.26"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.61%  "
$ws.Range("E36").Value = "  -6.39%  "
$ws.Range("E37").Value = "  -3.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.80"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "311.31"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -8.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.855"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -8.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.76"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.71"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.07%  "
$ws.Range("E43").Value = "  -0.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.70"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.581"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0936"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0521"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.03"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0229"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.48"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -9.51%  "
$ws.Range("D51").Value = "1.982.92"
